$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Replace the "Read PlcPcd" / "Calculate the largest..." /
#    "Grab all the characters..." block (paragraphs 11-13) with the
#    new, expanded 13-paragraph walkthrough.
# ------------------------------------------------------------------

$p11 = $d.Paragraphs(11)
$p13 = $d.Paragraphs(13)
$blockRange = $d.Range($p11.Range.Start, $p13.Range.End)
$blockRange.Delete()

$anchor = $d.Paragraphs(10).Range

$anchor.InsertParagraphAfter()
$p = $d.Paragraphs(11)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("Use lcb to calculate number of Pcds (lcb-4)/4+pcdSize). pcdSize is 8 (see section 2.8.35)")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(12)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("Add 1 to number of Pcds to get number of cps")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(13)
$p.Range.ListFormat.ListLevelNumber = 1
$p.Range.InsertAfter("Read PlcPcd")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(14)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("Read number of cps into aCp (4-byte unsigned integers)")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(15)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("Read number of Pcds into aPcd")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(16)
$p.Range.ListFormat.ListLevelNumber = 1
$p.Range.InsertAfter("Read Pcd")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(17)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("Read bytes 2-5 into FcCompressed structure")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(18)
$p.Range.ListFormat.ListLevelNumber = 1
$p.Range.InsertAfter("Read FcCompressed")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(19)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("Read 4 bytes as little-endian byte order and little-endian bit order values")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(20)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("First 30 bits (lowest order bits) are the FcCompressed. This includes all of bytes 0-2 and the first 2 bits in byte 3. Read as an unsigned integer")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(21)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("Read bit 31 (byte 3, bit 3) into fCompressed. 0 means it isn" + [char]0x2019 + "t compressed, 1 means it is")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(22)
$p.Range.ListFormat.ListLevelNumber = 1
$p.Range.InsertAfter("Loop through aPcd and Read characters from WordDocument stream")

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(23)
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.InsertAfter("Seek to aPcd[i].fc.fc")

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the "Prl" run-split to the end
#    of the text we just typed. Re-adding a bookmark with the same
#    name moves it (and the old split is no longer needed once the
#    bookmark that forced it is gone, so merge "P"+"rl" back to "Prl").
# ------------------------------------------------------------------

$goBackRange = $p.Range.Duplicate
$goBackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRange)

$d.Content.Find.Execute("Prl", $true, $false, $false, $false, $false, $true, 1, $false, "Prl", 2) | Out-Null
